$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '42.767.92'
Set-TextValue 'E2' '  -0.08%  '
Set-TextValue 'D3' '2.312.06'
Set-TextValue 'E3' '  +0.32%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '312.34'
Set-TextValue 'E5' '  -1.95%  '
Set-TextValue 'D6' '107.12'
Set-TextValue 'E6' '  +2.96%  '
Set-TextValue 'D7' '0.622'
Set-TextValue 'E7' '  -1.10%  '
Set-TextValue 'E8' '  -0.13%  '
Set-TextValue 'D9' '0.606'
Set-TextValue 'E9' '  +0.31%  '
Set-TextValue 'D10' '40.20'
Set-TextValue 'E10' '  +0.95%  '
Set-TextValue 'D11' '0.0912'
Set-TextValue 'E11' '  +0.30%  '
Set-TextValue 'E12' '  -1.75%  '
Set-TextValue 'E13' '  -1.67%  '
Set-TextValue 'E14' '  -1.65%  '
Set-TextValue 'D15' '15.35'
Set-TextValue 'E15' '  -0.33%  '
Set-TextValue 'D16' '2.663.15'
Set-TextValue 'E16' '  +0.34%  '
Set-TextValue 'D17' '2.307.74'
Set-TextValue 'E17' '  +0.25%  '
Set-TextValue 'D18' '42.930.31'
Set-TextValue 'E18' '  +0.76%  '
Set-TextValue 'E19' '  -0.90%  '
Set-TextValue 'E20' '  -0.41%  '
Set-TextValue 'D21' '13.16'
Set-TextValue 'E21' '  -10.35%  '
Set-TextValue 'D22' '73.60'
Set-TextValue 'E22' '  -0.48%  '
Set-TextValue 'E23' '  -1.85%  '
Set-TextValue 'D24' '266.22'
Set-TextValue 'E24' '  -0.32%  '
Set-TextValue 'D25' '2.25'
Set-TextValue 'E25' '  +0.50%  '
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'D27' '7.78'
Set-TextValue 'E27' '  +16.10%  '
Set-TextValue 'D28' '11.02'
Set-TextValue 'E28' '  +0.37%  '
Set-TextValue 'D29' '2.31'
Set-TextValue 'E29' '  +1.39%  '
Set-TextValue 'D30' '38.66'
Set-TextValue 'E30' '  +3.31%  '
Set-TextValue 'D31' '22.44'
Set-TextValue 'E31' '  -0.87%  '
Set-TextValue 'D32' '166.12'
Set-TextValue 'E32' '  +0.07%  '
Set-TextValue 'D33' '0.0873'
Set-TextValue 'E33' '  -1.21%  '
Set-TextValue 'D34' '2.75'
Set-TextValue 'E34' '  +5.15%  '
Set-TextValue 'E36' '  -1.02%  '
Set-TextValue 'E37' '  +1.56%  '
Set-TextValue 'E38' '  +0.89%  '
Set-TextValue 'E39' '  +4.61%  '
Set-TextValue 'D40' '3.68'
Set-TextValue 'E40' '  -1.28%  '
Set-TextValue 'D41' '1.60'
Set-TextValue 'E41' '  +0.50%  '
Set-TextValue 'D42' '103.39'
Set-TextValue 'E42' '  +7.91%  '
Set-TextValue 'D43' '71.03'
Set-TextValue 'E43' '  +0.73%  '
Set-TextValue 'D44' '0.234'
Set-TextValue 'E44' '  +1.86%  '
Set-TextValue 'D45' '1.00'
Set-TextValue 'E45' '  +0.44%  '
Set-TextValue 'D46' '12.74'
Set-TextValue 'E46' '  +3.97%  '
Set-TextValue 'D47' '112.85'
Set-TextValue 'E47' '  -1.87%  '
Set-TextValue 'B48' 'ordi'
Set-TextValue 'C48' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 'D48' '77.18'
Set-TextValue 'E48' '  -5.49%  '
Set-TextValue 'B49' 'Maker'
Set-TextValue 'C49' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D49' '1.656.35'
Set-TextValue 'E49' '  -2.10%  '
Set-TextValue 'D50' '8.82'
Set-TextValue 'E50' '  -0.14%  '
Set-TextValue 'D51' '5.25'
Set-TextValue 'E51' '  +1.07%  '
